# Clean up duplicate city rows in the Cidades list and refresh the AutoFilter.
#
# Three rows are duplicate / garbage entries that need to be removed:
#   - Row 12 "MARACAS"     is a duplicate of row 11 "MARACÁS"      -> merge totals, delete row 12
#   - Row 60 "RUI BARBOSA" is a duplicate of row 59 "RUY BARBOSA"  -> merge totals, delete row 60
#   - Row 96 "."           is a bogus entry with no real match     -> delete row 96 outright
#
# Deletions are performed bottom-to-top so earlier row numbers stay valid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the bogus "." row (row 96) — no merge target. ---
$ws.Rows(96).Delete()

# --- 2. Merge "RUI BARBOSA" (row 60) into "RUY BARBOSA" (row 59), then delete row 60. ---
$ws.Range("E59").Value = $ws.Range("E59").Value2 + $ws.Range("E60").Value2
$ws.Rows(60).Delete()

# --- 3. Merge "MARACAS" (row 12) into "MARACÁS" (row 11), then delete row 12. ---
$ws.Range("E11").Value = $ws.Range("E11").Value2 + $ws.Range("E12").Value2
$ws.Rows(12).Delete()

# --- 4. Re-apply the AutoFilter over the shrunk range A1:E120. ---
$rng = $ws.Range("A1:E120")
$rng.AutoFilter() | Out-Null
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$120")
$name.Visible = $false
